$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "doing my fi es"
$ws.Range("E3").Select()
